# Update the "想去人数" (want-to-go count) figures in the 展览 and 全部类型 sheets
# to reflect the latest generated output (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 805
$ws1.Range("F6").Value = 109
$ws1.Range("F8").Value = 4402
$ws1.Range("F10").Value = 4965
$ws1.Range("F11").Value = 556
$ws1.Range("F12").Value = 1250

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 805
$ws4.Range("F6").Value = 109
$ws4.Range("F9").Value = 4402
$ws4.Range("F11").Value = 4965
$ws4.Range("F12").Value = 556
$ws4.Range("F13").Value = 1250

$wb.Save()
